# Applies the "New links and commands" commit to the MQTT links workbook:
#   - 10 new rows (35-44) of topic/link reference rows appended to Hoja1
#   - rows 25-34 gain a (blank) column-C cell, matching the row-40 layout
#     that already had a 3rd link in column C
#   - column A-D widths were trimmed slightly
#   - row 3 height + the view scroll/selection moved to track the new rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 35-44: links/commands appended at the bottom of the list ---
$ws.Range("A35").Value = "Commandos Ethernet/Network (Linux)"
$ws.Range("B35").Value = "http://developer.toradex.com/knowledge-base/ethernet-network-(linux)"

$ws.Range("A36").Value = "Installing and Running applications on the Galileo Board Gen 2 "
$ws.Range("B36").Value = "https://docs.oracle.com/javame/8.3/get-started-galileo/installing-and-running-applications-intel-galileo-gen2-board.htm#MEEGG-GUID-F8F98050-3BE7-4E17-BFE5-72CCE9C4FBD1"

$ws.Range("A37").Value = "Configuración más completa Mosquitto Server"
$ws.Range("B37").Value = "http://lukse.lt/uzrasai/2015-02-internet-of-things-messaging-mqtt-1-installing-mosquitto-server/"

$ws.Range("A38").Value = "COMANDOS LINUX"
$ws.Range("B38").Value = "http://linuxcommand.org/man_pages/ls1.html"

$ws.Range("A39").Value = "Mosquito Guide (Broker,_sub,…..)"
$ws.Range("B39").Value = "https://mosquitto.org/documentation/"

$ws.Range("A40").Value = "SIMULADORES MODBUS"
$ws.Range("B40").Value = "http://docklight.de/"
$ws.Range("C40").Value = "http://www.plcsimulator.org/"

$ws.Range("A41").Value = "Tutorial de Modbus for Galileo"
$ws.Range("B41").Value = "https://www.cooking-hacks.com/documentation/tutorials/modbus-module-shield-tutorial-for-arduino-raspberry-pi-intel-galileo/"

$ws.Range("A42").Value = "Foro sobre Modbus y su implementación"
$ws.Range("B42").Value = "http://stackoverflow.com/questions/tagged/modbus?page=1&sort=newest&pagesize=15"

$ws.Range("A43").Value = "microcontroladores"
$ws.Range("B43").Value = "http://www.freertos.org/a00090.html"

$ws.Range("A44").Value = "nModbus (Otra libreria)"
$ws.Range("B44").Value = "http://www.mesta-automation.com/modbus-with-c-sharp-libraries-examples/"

# --- Rows 25-34: add the (empty, default-styled) column C cell that the
#     rest of the table already has, by cloning the blank-style C1 cell ---
$ws.Range("C1").Copy()
$ws.Range("C25:C34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 3 height tightened to match the rest of the data rows ---
$ws.Rows.Item(3).RowHeight = 13.8

# --- Column width tweaks (A-D narrowed slightly) ---
$ws.Columns.Item(1).ColumnWidth = 58.5
$ws.Columns.Item(2).ColumnWidth = 119.16666666666667
$ws.Columns.Item(3).ColumnWidth = 94.83333333333333
$ws.Columns.Item(4).ColumnWidth = 87.83333333333333

# --- Selection / scroll position follows the newly added rows ---
$ws.Range("A45").Select() | Out-Null
